$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the styles we need to re-apply later (style objects / names get
# reset by Hyperlinks.Add, so capture them up-front from cells that already
# carry the right formatting).
$hyperlinkStyle = $ws.Range("D4").Style

# --- Remove the "Fox" asset entry that used to live in row 3 -------------
# A3/B3 (type + where) are fully cleared (no cell left behind at all),
# while C3/D3 keep their existing fill/hyperlink formatting but lose their
# value - they become blank, styled placeholder cells.
$ws.Range("A3:B3").ClearContents()
$ws.Range("C3:D3").ClearContents()

# The old D3 hyperlink (pointing at the Fox asset page) must disappear.
# The engine only supports wiping *all* hyperlinks on the sheet at once, so
# remove them all and then recreate the two that should survive (D4 + D11)
# pointing at their original targets.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D4"), "https://assetstore.unity.com/packages/3d/props/simple-gems-ultimate-animated-customizable-pack-73764") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "https://free3d.com/3d-model/coin-4532.html") | Out-Null

# Adding a hyperlink resets cell formatting to the default hyperlink style -
# restore the original look (font/fill) that these cells had before.
$ws.Range("D4").Style = $hyperlinkStyle
$ws.Range("D11").Style = $hyperlinkStyle

# --- New empty (but styled) row 5, under the AurynSky/truphy row ---------
$ws.Range("D5").Style = $hyperlinkStyle

# --- Update the selected range shown when the sheet is opened ------------
[void]$ws.Range("C5:L5").Select()
